$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.168.14"
$ws.Range("E2").Value = "  -1.18%  "

$ws.Range("D3").Value = "1.659.06"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.27"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5202"
$ws.Range("E6").Value = "  -2.63%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2631"
$ws.Range("E8").Value = "  -3.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06267"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.77"
$ws.Range("E10").Value = "  -5.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07725"
$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.420"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.627.78"
$ws.Range("E13").Value = "  -3.32%  "

$ws.Range("D14").Value = "1.883.77"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5431"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").Value = "0.0₅8130"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.43"
$ws.Range("E17").Value = "  -2.11%  "

$ws.Range("D18").Value = "26.196.80"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.627"
$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.93"
$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  -2.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.058"

$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.74"
$ws.Range("E25").Value = "  -1.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1226"
$ws.Range("E26").Value = "  -4.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.181"
$ws.Range("E27").Value = "  -3.28%  "

$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.405"
$ws.Range("E29").Value = "  -2.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05975"
$ws.Range("E30").Value = "  -5.23%  "

$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.551"
$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.240"
$ws.Range("E33").Value = "  -6.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.612"
$ws.Range("E34").Value = "  -5.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9672"
$ws.Range("E35").Value = "  -4.48%  "

$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5681"
$ws.Range("E38").Value = "  -7.67%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.012"
$ws.Range("E39").Value = "  -2.53%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01597"
$ws.Range("E40").Value = "  -2.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8580"
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("D43").Value = "1.013.54"
$ws.Range("E43").Value = "  -7.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.41"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").Value = "1.799.74"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.80"
$ws.Range("E47").Value = "  -3.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.016"
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.455"
$ws.Range("E51").Value = "  -1.19%  "
